# Generate Report for handback
#
# The localization report workbook records, per-language, the handoff/
# handback lifecycle of each source file. This run marks the two tracked
# files as handed back ("in sync with en-US") and records their handback
# target/file columns + handback timestamp. The "Ready for handoff"
# status text (shared by the Overview roll-up sheet and each per-language
# sheet) is replaced everywhere with the new status text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status text is shared with the per-language sheets,
# so it must be refreshed here too (B2:C2, B3:C3) even though no other
# cell on this sheet changes.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- Per-language sheets: zh-cn handed back at 2016-01-26 07:37:20,
# de-de handed back at 2016-01-26 07:37:39.
$langs = @(
    @{ Sheet = "zh-cn"; HandbackTime = "2016-01-26 07:37:20" },
    @{ Sheet = "de-de"; HandbackTime = "2016-01-26 07:37:39" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Row 2 and row 3 are the two real source files; row 4 is the
    # .localization-config row, which stays "Ignored" and is untouched.
    foreach ($row in 2, 3) {
        $ws.Range("B$row").Value = $newStatus
    }

    # Capture the existing hyperlink addresses/text for the "source md"
    # (column A) and "handoff target file" (column C) links on row 2 so
    # the new "Latest Target File" (E) / "Latest Handback File" (F)
    # columns can mirror them exactly, as the report generator does.
    $aLink = $null
    $cLink = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$A$2') { $aLink = $h }
        if ($h.Range.Address() -eq '$C$2') { $cLink = $h }
    }

    foreach ($row in 2, 3) {
        $eCell = $ws.Range("E$row")
        $eCell.Value = $aLink.TextToDisplay()
        $ws.Hyperlinks.Add($eCell, $aLink.Address(), "", "", $aLink.TextToDisplay())

        $fCell = $ws.Range("F$row")
        $fCell.Value = $cLink.TextToDisplay()
        $ws.Hyperlinks.Add($fCell, $cLink.Address(), "", "", $cLink.TextToDisplay())

        # Latest Handback DateTime (column G): was the zero-date sentinel,
        # now the real handback timestamp for this language.
        $ws.Range("G$row").Value = $lang.HandbackTime
    }
}
